$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 16.83279933333333
$ws.Cells.Item(2, 8).Value = 50.498398
$ws.Cells.Item(2, 9).Value = 0.04383102208811961
$ws.Cells.Item(2, 10).Value = 0.04383102208811961
$ws.Cells.Item(2, 13).Value = 2.598166333333333
$ws.Cells.Item(2, 14).Value = 7.794499
$ws.Cells.Item(2, 15).Value = 0.3466013321552429
$ws.Cells.Item(2, 16).Value = 0.3466013321552429
$ws.Cells.Item(2, 17).Value = 43.73441252362245
$ws.Cells.Item(2, 18).Value = 393.609712712602
$ws.Cells.Item(2, 19).Value = 0.01519189064546814
$ws.Cells.Item(2, 20).Value = 0.01519189064546814
$ws.Cells.Item(3, 7).Value = 16.83279933333333
$ws.Cells.Item(3, 8).Value = 50.498398
$ws.Cells.Item(3, 9).Value = 0.04383102208811961
$ws.Cells.Item(3, 10).Value = 0.04383102208811961
$ws.Cells.Item(3, 15).Value = 0.5780859172985858
$ws.Cells.Item(3, 16).Value = 0.5780859172985858
$ws.Cells.Item(3, 17).Value = 72.94330874039778
$ws.Cells.Item(3, 18).Value = 656.48977866358
$ws.Cells.Item(3, 19).Value = 0.0253380966099452
$ws.Cells.Item(3, 20).Value = 0.0253380966099452
$ws.Cells.Item(4, 7).Value = 16.83279933333333
$ws.Cells.Item(4, 8).Value = 50.498398
$ws.Cells.Item(4, 9).Value = 0.04383102208811961
$ws.Cells.Item(4, 10).Value = 0.04383102208811961
$ws.Cells.Item(4, 13).Value = 0.4692043333333333
$ws.Cells.Item(4, 14).Value = 1.407613
$ws.Cells.Item(4, 15).Value = 0.06259293136852516
$ws.Cells.Item(4, 16).Value = 0.06259293136852516
$ws.Cells.Item(4, 17).Value = 7.898022389330444
$ws.Cells.Item(4, 18).Value = 71.082201503974
$ws.Cells.Item(4, 19).Value = 0.002743512157373981
$ws.Cells.Item(4, 20).Value = 0.002743512157373981
$ws.Cells.Item(5, 7).Value = 16.83279933333333
$ws.Cells.Item(5, 8).Value = 50.498398
$ws.Cells.Item(5, 9).Value = 0.04383102208811961
$ws.Cells.Item(5, 10).Value = 0.04383102208811961
$ws.Cells.Item(5, 11).Value = 2.0
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.09534933333333333
$ws.Cells.Item(5, 14).Value = 0.286048
$ws.Cells.Item(5, 15).Value = 0.01271981917764605
$ws.Cells.Item(5, 16).Value = 0.01271981917764604
$ws.Cells.Item(5, 17).Value = 1.604996194567111
$ws.Cells.Item(5, 18).Value = 14.444965751104
$ws.Cells.Item(5, 19).Value = 0.0005575226753322912
$ws.Cells.Item(5, 20).Value = 0.0005575226753322912
$ws.Cells.Item(6, 9).Value = 0.8654671168650652
$ws.Cells.Item(6, 10).Value = 0.8654671168650654
$ws.Cells.Item(6, 13).Value = 2.598166333333333
$ws.Cells.Item(6, 14).Value = 7.794499
$ws.Cells.Item(6, 15).Value = 0.3466013321552429
$ws.Cells.Item(6, 16).Value = 0.3466013321552429
$ws.Cells.Item(6, 17).Value = 863.5595090278842
$ws.Cells.Item(6, 18).Value = 7772.035581250959
$ws.Cells.Item(6, 19).Value = 0.2999720556419889
$ws.Cells.Item(6, 20).Value = 0.299972055641989
$ws.Cells.Item(7, 9).Value = 0.8654671168650652
$ws.Cells.Item(7, 10).Value = 0.8654671168650654
$ws.Cells.Item(7, 15).Value = 0.5780859172985858
$ws.Cells.Item(7, 16).Value = 0.5780859172985858
$ws.Cells.Item(7, 19).Value = 0.5003143521447035
$ws.Cells.Item(7, 20).Value = 0.5003143521447037
$ws.Cells.Item(8, 9).Value = 0.8654671168650652
$ws.Cells.Item(8, 10).Value = 0.8654671168650654
$ws.Cells.Item(8, 13).Value = 0.4692043333333333
$ws.Cells.Item(8, 14).Value = 1.407613
$ws.Cells.Item(8, 15).Value = 0.06259293136852516
$ws.Cells.Item(8, 16).Value = 0.06259293136852516
$ws.Cells.Item(8, 17).Value = 155.9507020504162
$ws.Cells.Item(8, 18).Value = 1403.556318453746
$ws.Cells.Item(8, 19).Value = 0.05417212384765037
$ws.Cells.Item(8, 20).Value = 0.05417212384765038
$ws.Cells.Item(9, 9).Value = 0.8654671168650652
$ws.Cells.Item(9, 10).Value = 0.8654671168650654
$ws.Cells.Item(9, 11).Value = 2.0
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.09534933333333333
$ws.Cells.Item(9, 14).Value = 0.286048
$ws.Cells.Item(9, 15).Value = 0.01271981917764605
$ws.Cells.Item(9, 16).Value = 0.01271981917764604
$ws.Cells.Item(9, 17).Value = 31.69151351977955
$ws.Cells.Item(9, 18).Value = 285.223621678016
$ws.Cells.Item(9, 19).Value = 0.01100858523072229
$ws.Cells.Item(9, 20).Value = 0.01100858523072229
$ws.Cells.Item(10, 7).Value = 34.50825133333333
$ws.Cells.Item(10, 8).Value = 103.524754
$ws.Cells.Item(10, 9).Value = 0.08985623225594501
$ws.Cells.Item(10, 10).Value = 0.08985623225594502
$ws.Cells.Item(10, 13).Value = 2.598166333333333
$ws.Cells.Item(10, 14).Value = 7.794499
$ws.Cells.Item(10, 15).Value = 0.3466013321552429
$ws.Cells.Item(10, 16).Value = 0.3466013321552429
$ws.Cells.Item(10, 17).Value = 89.65817683647178
$ws.Cells.Item(10, 18).Value = 806.923591528246
$ws.Cells.Item(10, 19).Value = 0.03114428980236145
$ws.Cells.Item(10, 20).Value = 0.03114428980236145
$ws.Cells.Item(11, 7).Value = 34.50825133333333
$ws.Cells.Item(11, 8).Value = 103.524754
$ws.Cells.Item(11, 9).Value = 0.08985623225594501
$ws.Cells.Item(11, 10).Value = 0.08985623225594502
$ws.Cells.Item(11, 15).Value = 0.5780859172985858
$ws.Cells.Item(11, 16).Value = 0.5780859172985858
$ws.Cells.Item(11, 17).Value = 149.5381713553711
$ws.Cells.Item(11, 18).Value = 1345.84354219834
$ws.Cells.Item(11, 19).Value = 0.05194462244867275
$ws.Cells.Item(11, 20).Value = 0.05194462244867275
$ws.Cells.Item(12, 7).Value = 34.50825133333333
$ws.Cells.Item(12, 8).Value = 103.524754
$ws.Cells.Item(12, 9).Value = 0.08985623225594501
$ws.Cells.Item(12, 10).Value = 0.08985623225594502
$ws.Cells.Item(12, 13).Value = 0.4692043333333333
$ws.Cells.Item(12, 14).Value = 1.407613
$ws.Cells.Item(12, 15).Value = 0.06259293136852516
$ws.Cells.Item(12, 16).Value = 0.06259293136852516
$ws.Cells.Item(12, 17).Value = 16.19142106135578
$ws.Cells.Item(12, 18).Value = 145.722789552202
$ws.Cells.Item(12, 19).Value = 0.005624364978630622
$ws.Cells.Item(12, 20).Value = 0.005624364978630623
$ws.Cells.Item(13, 7).Value = 34.50825133333333
$ws.Cells.Item(13, 8).Value = 103.524754
$ws.Cells.Item(13, 9).Value = 0.08985623225594501
$ws.Cells.Item(13, 10).Value = 0.08985623225594502
$ws.Cells.Item(13, 11).Value = 2.0
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.09534933333333333
$ws.Cells.Item(13, 14).Value = 0.286048
$ws.Cells.Item(13, 15).Value = 0.01271981917764605
$ws.Cells.Item(13, 16).Value = 0.01271981917764604
$ws.Cells.Item(13, 17).Value = 3.290338759132444
$ws.Cells.Item(13, 18).Value = 29.613048832192
$ws.Cells.Item(13, 19).Value = 0.001142955026280187
$ws.Cells.Item(13, 20).Value = 0.001142955026280187
$ws.Cells.Item(14, 7).Value = 0.324754
$ws.Cells.Item(14, 8).Value = 0.974262
$ws.Cells.Item(14, 9).Value = 0.0008456287908700705
$ws.Cells.Item(14, 10).Value = 0.0008456287908700706
$ws.Cells.Item(14, 13).Value = 2.598166333333333
$ws.Cells.Item(14, 14).Value = 7.794499
$ws.Cells.Item(14, 15).Value = 0.3466013321552429
$ws.Cells.Item(14, 16).Value = 0.3466013321552429
$ws.Cells.Item(14, 17).Value = 0.8437649094153333
$ws.Cells.Item(14, 18).Value = 7.593884184738
$ws.Cells.Item(14, 19).Value = 0.0002930960654243937
$ws.Cells.Item(14, 20).Value = 0.0002930960654243938
$ws.Cells.Item(15, 7).Value = 0.324754
$ws.Cells.Item(15, 8).Value = 0.974262
$ws.Cells.Item(15, 9).Value = 0.0008456287908700705
$ws.Cells.Item(15, 10).Value = 0.0008456287908700706
$ws.Cells.Item(15, 15).Value = 0.5780859172985858
$ws.Cells.Item(15, 16).Value = 0.5780859172985858
$ws.Cells.Item(15, 17).Value = 1.407290066113333
$ws.Cells.Item(15, 18).Value = 12.66561059502
$ws.Cells.Item(15, 19).Value = 0.0004888460952642186
$ws.Cells.Item(15, 20).Value = 0.0004888460952642187
$ws.Cells.Item(16, 7).Value = 0.324754
$ws.Cells.Item(16, 8).Value = 0.974262
$ws.Cells.Item(16, 9).Value = 0.0008456287908700705
$ws.Cells.Item(16, 10).Value = 0.0008456287908700706
$ws.Cells.Item(16, 13).Value = 0.4692043333333333
$ws.Cells.Item(16, 14).Value = 1.407613
$ws.Cells.Item(16, 15).Value = 0.06259293136852516
$ws.Cells.Item(16, 16).Value = 0.06259293136852516
$ws.Cells.Item(16, 17).Value = 0.1523759840673333
$ws.Cells.Item(16, 18).Value = 1.371383856606
$ws.Cells.Item(16, 19).Value = 0.00005293038487017923
$ws.Cells.Item(16, 20).Value = 0.00005293038487017924
$ws.Cells.Item(17, 7).Value = 0.324754
$ws.Cells.Item(17, 8).Value = 0.974262
$ws.Cells.Item(17, 9).Value = 0.0008456287908700705
$ws.Cells.Item(17, 10).Value = 0.0008456287908700706
$ws.Cells.Item(17, 11).Value = 2.0
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.09534933333333333
$ws.Cells.Item(17, 14).Value = 0.286048
$ws.Cells.Item(17, 15).Value = 0.01271981917764605
$ws.Cells.Item(17, 16).Value = 0.01271981917764604
$ws.Cells.Item(17, 17).Value = 0.03096507739733333
$ws.Cells.Item(17, 18).Value = 0.278685696576
$ws.Cells.Item(17, 19).Value = 0.00001075624531127876
$ws.Cells.Item(17, 20).Value = 0.00001075624531127876
